$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 41
$ws.Range("B2").Value = 153
$ws.Range("B3").Value = 247
$ws.Range("B4").Value = 304
